# New test case added for spreadsheet importing.
#
# This reproduces, via Excel COM-interop calls, the changes made to the
# "Sheet9" worksheet (a new "fourth" column header + a date/time value in a
# new position, with the old column-A values shifted one column to the
# right) together with the supporting column-width/number-format tweaks.
# The "testAccounts" sheet and the shared-strings table pick up their index
# shift automatically because the new "fourth" string is interned in the
# same relative position as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet9")

# --- New header cell: D1 = "fourth" ------------------------------------
$ws.Range("D1").Value = "fourth"

# --- Row 2: shift "a"/45 right one column, add a new timestamp in A2 ---
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = 45
$ws.Range("A2").Value = 40909.5163657407
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd\ hh:mm"

# --- Row 3: shift "s"/5 right one column, clear the old A3 -------------
$ws.Range("B3").Value = "s"
$ws.Range("C3").Value = 5
$ws.Range("A3").ClearContents()

# --- Row 4: shift "d"/67 right one column, clear the old A4 ------------
$ws.Range("B4").Value = "d"
$ws.Range("C4").Value = 67
$ws.Range("A4").ClearContents()

# --- Column widths (characters) matching the source edit ---------------
$ws.Columns.Item(1).ColumnWidth = 15.97
$ws.Columns.Item(3).ColumnWidth = 16.53

# --- Selection left on D13 after the edit -------------------------------
$ws.Activate()
$ws.Range("D13").Select()
